$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11993.5
$ws.Range("I43").Value = 9889.700000000001
$ws.Range("K43").Value = 9889.700000000001
$ws.Range("M43").Value = -9820.700000000001

$ws.Range("H76").Value = 5750
$ws.Range("J76").Value = 8000
$ws.Range("L76").Value = 8000
$ws.Range("N76").Value = -8630

$ws.Range("H79").Value = 5750
$ws.Range("J79").Value = 8000
$ws.Range("L79").Value = 8000
$ws.Range("N79").Value = -10184

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H132").Value = 844.3871
$ws.Range("I132").Value = 861.4828
$ws.Range("K132").Value = 2584.4484
$ws.Range("M132").Value = -54.44840000000022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 6950
$ws.Range("J14").Value = 6950
$ws.Range("L14").Value = 6950
$ws.Range("N14").Value = -7300

$ws.Range("H46").Value = 24276.715
$ws.Range("I46").Value = 27401
$ws.Range("J46").Value = 20111
$ws.Range("K46").Value = 27401
$ws.Range("L46").Value = 20111
$ws.Range("M46").Value = -27082
$ws.Range("N46").Value = -20749

$ws.Range("H74").Value = 758.3333
$ws.Range("I74").Value = 758.3333
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 758.3333
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 115.6667
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 758.3333
$ws.Range("I77").Value = 758.3333
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3791.6665
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 576.3334999999997
$ws.Range("N77").ClearContents()

$ws.Range("H88").Value = 1255.1111
$ws.Range("I88").Value = 1032.25
$ws.Range("J88").Value = 1433.4
$ws.Range("K88").Value = 1032.25
$ws.Range("L88").Value = 1433.4
$ws.Range("M88").Value = -626.25
$ws.Range("N88").Value = -2245.4

$ws.Range("H91").Value = 1255.1111
$ws.Range("I91").Value = 1032.25
$ws.Range("J91").Value = 1433.4
$ws.Range("K91").Value = 1032.25
$ws.Range("L91").Value = 1433.4
$ws.Range("M91").Value = 371.75
$ws.Range("N91").Value = -4241.4

$ws.Range("H104").Value = 54000
$ws.Range("J104").Value = 54000
$ws.Range("L104").Value = 54000
$ws.Range("N104").Value = -60988

$ws.Range("H110").Value = 2679
$ws.Range("I110").Value = 2222.0625
$ws.Range("K110").Value = 2222.0625
$ws.Range("M110").Value = -177.0625

$ws.Range("H122").Value = 9999999
$ws.Range("I122").Value = 9999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 29999997
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -29997547
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 18590.223
$ws.Range("I132").Value = 2308.8
$ws.Range("J132").Value = 99997.336
$ws.Range("K132").Value = 6926.400000000001
$ws.Range("L132").Value = 299992.008
$ws.Range("M132").Value = -4396.400000000001
$ws.Range("N132").Value = -305052.008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4330.136
$ws.Range("I31").Value = 3455.1875
$ws.Range("K31").Value = 3455.1875
$ws.Range("M31").Value = -3160.1875

$ws.Range("H34").Value = 4330.136
$ws.Range("I34").Value = 3455.1875
$ws.Range("K34").Value = 3455.1875
$ws.Range("M34").Value = -3253.1875

$ws.Range("H58").Value = 1854.5641
$ws.Range("I58").Value = 1085.3667
$ws.Range("J58").Value = 4418.5557
$ws.Range("K58").Value = 1085.3667
$ws.Range("L58").Value = 4418.5557
$ws.Range("M58").Value = -882.3667
$ws.Range("N58").Value = -4824.5557

$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992

$ws.Range("H134").Value = 2322.72
$ws.Range("I134").Value = 2191.0278
$ws.Range("J134").Value = 2661.3572
$ws.Range("K134").Value = 6573.0834
$ws.Range("L134").Value = 7984.071599999999
$ws.Range("M134").Value = -4038.0834
$ws.Range("N134").Value = -13054.0716

$ws.Range("H136").Value = 1854.5641
$ws.Range("I136").Value = 1085.3667
$ws.Range("J136").Value = 4418.5557
$ws.Range("K136").Value = 3256.1001
$ws.Range("L136").Value = 13255.6671
$ws.Range("M136").Value = -706.1001000000001
$ws.Range("N136").Value = -18355.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 7666.6665
$ws.Range("J138").Value = 7666.6665
$ws.Range("L138").Value = 22999.9995
$ws.Range("N138").Value = -33279.99950000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 26948
$ws.Range("J57").Value = 21716.2
$ws.Range("L57").Value = 21716.2
$ws.Range("N57").Value = -23356.2

$ws.Range("H126").Value = 3500.4443
$ws.Range("I126").Value = 3572.2856
$ws.Range("K126").Value = 10716.8568
$ws.Range("M126").Value = -8246.856800000001

$ws.Range("H141").Value = 30000
$ws.Range("J141").Value = 30000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 13998.333
$ws.Range("J3").Value = 13998.333
$ws.Range("L3").Value = 13998.333
$ws.Range("N3").Value = -14222.333

$ws.Range("H7").Value = 2131.3333
$ws.Range("I7").Value = 2131.3333
$ws.Range("K7").Value = 2131.3333
$ws.Range("M7").Value = -2019.3333

$ws.Range("H15").Value = 13998.333
$ws.Range("J15").Value = 13998.333
$ws.Range("L15").Value = 13998.333
$ws.Range("N15").Value = -14338.333

$ws.Range("H43").Value = 8072.4546
$ws.Range("J43").Value = 8079.8
$ws.Range("L43").Value = 8079.8
$ws.Range("N43").Value = -8465.799999999999

$ws.Range("H126").Value = 2131.3333
$ws.Range("I126").Value = 2131.3333
$ws.Range("K126").Value = 6393.999899999999
$ws.Range("M126").Value = -3923.999899999999

$ws.Range("H136").Value = 7176.2666
$ws.Range("I136").Value = 7096.6924
$ws.Range("K136").Value = 21290.0772
$ws.Range("M136").Value = -18740.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 1875
$ws.Range("I51").Value = 1875
$ws.Range("K51").Value = 1875
$ws.Range("M51").Value = -1365

$ws.Range("H101").Value = 3694.25
$ws.Range("J101").Value = 3694.25
$ws.Range("L101").Value = 3694.25
$ws.Range("N101").Value = -10184.25

$ws.Range("H107").Value = 199.75
$ws.Range("I107").Value = 199.75
$ws.Range("K107").Value = 599.25
$ws.Range("M107").Value = 1320.75

$ws.Range("H132").Value = 3217.238
$ws.Range("J132").Value = 7250
$ws.Range("L132").Value = 21750
$ws.Range("N132").Value = -26810

$ws.Range("H136").Value = 1112.6666
$ws.Range("I136").Value = 751.7917
$ws.Range("K136").Value = 2255.3751
$ws.Range("M136").Value = 294.6248999999998
